$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4200
$ws.Range("I62").Value = 10000
$ws.Range("J62").Value = 1300
$ws.Range("K62").Value = 10000
$ws.Range("L62").Value = 1300
$ws.Range("M62").Value = -9376
$ws.Range("N62").Value = -2548

$ws.Range("H65").Value = 4200
$ws.Range("I65").Value = 10000
$ws.Range("J65").Value = 1300
$ws.Range("K65").Value = 50000
$ws.Range("L65").Value = 6500
$ws.Range("M65").Value = -46880
$ws.Range("N65").Value = -12740

$ws.Range("H76").Value = 5575.846
$ws.Range("I76").Value = 4467.381
$ws.Range("J76").Value = 6869.0557
$ws.Range("K76").Value = 4467.381
$ws.Range("L76").Value = 6869.0557
$ws.Range("M76").Value = -4152.381
$ws.Range("N76").Value = -7499.0557

$ws.Range("H79").Value = 5575.846
$ws.Range("I79").Value = 4467.381
$ws.Range("J79").Value = 6869.0557
$ws.Range("K79").Value = 4467.381
$ws.Range("L79").Value = 6869.0557
$ws.Range("M79").Value = -3375.381
$ws.Range("N79").Value = -9053.055700000001

$ws.Range("H86").Value = 3999.4546
$ws.Range("I86").Value = 10000
$ws.Range("J86").Value = 3399.4
$ws.Range("K86").Value = 10000
$ws.Range("L86").Value = 3399.4
$ws.Range("M86").Value = -8877
$ws.Range("N86").Value = -5645.4

$ws.Range("H89").Value = 3999.4546
$ws.Range("I89").Value = 10000
$ws.Range("J89").Value = 3399.4
$ws.Range("K89").Value = 50000
$ws.Range("L89").Value = 16997
$ws.Range("M89").Value = -44384
$ws.Range("N89").Value = -28229

$ws.Range("H125").Value = 2380.6365
$ws.Range("I125").Value = 4337.4
$ws.Range("J125").Value = 750
$ws.Range("K125").Value = 39036.6
$ws.Range("L125").Value = 6750
$ws.Range("M125").Value = -36576.6
$ws.Range("N125").Value = -11670

$ws.Range("H127").Value = 802.625
$ws.Range("I127").Value = 531.55554
$ws.Range("J127").Value = 908.6957
$ws.Range("K127").Value = 1594.66662
$ws.Range("L127").Value = 2726.0871
$ws.Range("M127").Value = 3365.33338
$ws.Range("N127").Value = -12646.0871

$ws.Range("H137").Value = 6061891
$ws.Range("I137").Value = 812.5833
$ws.Range("J137").Value = 22224766
$ws.Range("K137").Value = 2437.7499
$ws.Range("L137").Value = 66674298
$ws.Range("M137").Value = 112.2501000000002
$ws.Range("N137").Value = -66679398

$ws.Range("H138").Value = 1887.6097
$ws.Range("I138").Value = 1236.2307
$ws.Range("J138").Value = 3016.6667
$ws.Range("K138").Value = 3708.6921
$ws.Range("L138").Value = 9050.000100000001
$ws.Range("M138").Value = 1431.3079
$ws.Range("N138").Value = -19330.0001

$ws.Range("H139").Value = 41933.332
$ws.Range("J139").Value = 41933.332
$ws.Range("L139").Value = 41933.332
$ws.Range("N139").Value = -52213.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 13516823
$ws.Range("I61").Value = 18521728
$ws.Range("J61").Value = 3579.4
$ws.Range("K61").Value = 18521728
$ws.Range("L61").Value = 3579.4
$ws.Range("M61").Value = -18521516
$ws.Range("N61").Value = -4003.4

$ws.Range("H102").Value = 2500
$ws.Range("I102").Value = 1000
$ws.Range("J102").Value = 2875
$ws.Range("K102").Value = 1000
$ws.Range("L102").Value = 2875
$ws.Range("M102").Value = 622
$ws.Range("N102").Value = -6119

$ws.Range("H110").Value = 1000
$ws.Range("I110").Value = 1000
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1000
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1045
$ws.Range("N110").ClearContents()

$ws.Range("H122").Value = 8592.263000000001
$ws.Range("I122").Value = 10374.286
$ws.Range("J122").Value = 3602.6
$ws.Range("K122").Value = 31122.858
$ws.Range("L122").Value = 10807.8
$ws.Range("M122").Value = -28672.858
$ws.Range("N122").Value = -15707.8

$ws.Range("H132").Value = 7578045
$ws.Range("I132").Value = 15627191
$ws.Range("J132").Value = 2378.353
$ws.Range("K132").Value = 46881573
$ws.Range("L132").Value = 7135.059
$ws.Range("M132").Value = -46879043
$ws.Range("N132").Value = -12195.059

$ws.Range("H136").Value = 13516823
$ws.Range("I136").Value = 18521728
$ws.Range("J136").Value = 3579.4
$ws.Range("K136").Value = 55565184
$ws.Range("L136").Value = 10738.2
$ws.Range("M136").Value = -55562634
$ws.Range("N136").Value = -15838.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3724.4443
$ws.Range("I134").Value = 2395.6843
$ws.Range("J134").Value = 6880.25
$ws.Range("K134").Value = 7187.0529
$ws.Range("L134").Value = 20640.75
$ws.Range("M134").Value = -4652.0529
$ws.Range("N134").Value = -25710.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 648.5
$ws.Range("I22").Value = 540
$ws.Range("J22").Value = 757
$ws.Range("K22").Value = 540
$ws.Range("L22").Value = 757
$ws.Range("M22").Value = -190
$ws.Range("N22").Value = -1457

$ws.Range("H86").Value = 2326.923
$ws.Range("I86").Value = 2734.5
$ws.Range("J86").Value = 1977.5714
$ws.Range("K86").Value = 2734.5
$ws.Range("L86").Value = 1977.5714
$ws.Range("M86").Value = -1611.5
$ws.Range("N86").Value = -4223.5714

$ws.Range("H89").Value = 2326.923
$ws.Range("I89").Value = 2734.5
$ws.Range("J89").Value = 1977.5714
$ws.Range("K89").Value = 13672.5
$ws.Range("L89").Value = 9887.857
$ws.Range("M89").Value = -8056.5
$ws.Range("N89").Value = -21119.857

$ws.Range("H132").Value = 2501.1072
$ws.Range("I132").Value = 1856.75
$ws.Range("J132").Value = 3360.25
$ws.Range("K132").Value = 5570.25
$ws.Range("L132").Value = 10080.75
$ws.Range("M132").Value = -3040.25
$ws.Range("N132").Value = -15140.75

$ws.Range("H140").Value = 28153.334
$ws.Range("J140").Value = 28153.334
$ws.Range("L140").Value = 28153.334
$ws.Range("N140").Value = -38513.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3998
$ws.Range("I5").Value = 622
$ws.Range("K5").Value = 1866
$ws.Range("M5").Value = -1754

$ws.Range("H25").Value = 2700

$ws.Range("H30").Value = 2700

$ws.Range("H135").Value = 3998
$ws.Range("I135").Value = 622
$ws.Range("K135").Value = 5598
$ws.Range("M135").Value = -3063

$ws.Range("H137").Value = 9263899
$ws.Range("I137").Value = 18520316
$ws.Range("J137").Value = 7481.4443
$ws.Range("K137").Value = 55560948
$ws.Range("L137").Value = 22444.3329
$ws.Range("M137").Value = -55555848
$ws.Range("N137").Value = -32644.3329

$ws.Range("H138").Value = 2663.8572
$ws.Range("I138").Value = 1466.1666
$ws.Range("J138").Value = 9850
$ws.Range("K138").Value = 4398.4998
$ws.Range("L138").Value = 29550
$ws.Range("M138").Value = 741.5002000000004
$ws.Range("N138").Value = -39830

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2565530.8
$ws.Range("I122").Value = 3704848.5
$ws.Range("J122").Value = 2065.5
$ws.Range("K122").Value = 11114545.5
$ws.Range("L122").Value = 6196.5
$ws.Range("M122").Value = -11112095.5
$ws.Range("N122").Value = -11096.5

$ws.Range("H132").Value = 4160.3096
$ws.Range("I132").Value = 3687.087
$ws.Range("J132").Value = 4733.1577
$ws.Range("K132").Value = 11061.261
$ws.Range("L132").Value = 14199.4731
$ws.Range("M132").Value = -8531.261
$ws.Range("N132").Value = -19259.4731

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

$ws.Range("H40").Value = 3938.3872
$ws.Range("I40").Value = 6828
$ws.Range("K40").Value = 6828
$ws.Range("M40").Value = -6692

$ws.Range("H68").Value = 2902.625
$ws.Range("I68").Value = 3091.6667
$ws.Range("J68").Value = 2789.2
$ws.Range("K68").Value = 3091.6667
$ws.Range("L68").Value = 2789.2
$ws.Range("M68").Value = -2342.6667
$ws.Range("N68").Value = -4287.2

$ws.Range("H71").Value = 2902.625
$ws.Range("I71").Value = 3091.6667
$ws.Range("J71").Value = 2789.2
$ws.Range("K71").Value = 15458.3335
$ws.Range("L71").Value = 13946
$ws.Range("M71").Value = -11714.3335
$ws.Range("N71").Value = -21434

$ws.Range("H122").Value = 7127.8096
$ws.Range("I122").Value = 9559.875
$ws.Range("J122").Value = 5631.154
$ws.Range("K122").Value = 28679.625
$ws.Range("L122").Value = 16893.462
$ws.Range("M122").Value = -26229.625
$ws.Range("N122").Value = -21793.462

$ws.Range("H136").Value = 27788020
$ws.Range("I136").Value = 71430530
$ws.Range("J136").Value = 15515.272
$ws.Range("K136").Value = 214291590
$ws.Range("L136").Value = 46545.81600000001
$ws.Range("M136").Value = -214289040
$ws.Range("N136").Value = -51645.81600000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2415.8386
$ws.Range("J96").Value = 2847.75
$ws.Range("L96").Value = 2847.75
$ws.Range("N96").Value = -5593.75

$ws.Range("H122").Value = 2461.875
$ws.Range("I122").Value = 2288.889
$ws.Range("J122").Value = 2684.2856
$ws.Range("K122").Value = 6866.667
$ws.Range("L122").Value = 8052.8568
$ws.Range("M122").Value = -4416.667
$ws.Range("N122").Value = -12952.8568

$ws.Range("H125").Value = 60314
$ws.Range("J125").Value = 60314
$ws.Range("L125").Value = 60314
$ws.Range("N125").Value = -70154

$ws.Range("H126").Value = 2492.7144
$ws.Range("I126").Value = 1738.3529
$ws.Range("J126").Value = 5698.75
$ws.Range("K126").Value = 5215.0587
$ws.Range("L126").Value = 17096.25
$ws.Range("M126").Value = -2745.0587
$ws.Range("N126").Value = -22036.25

$ws.Range("H132").Value = 2141.2
$ws.Range("I132").Value = 1112.5714
$ws.Range("J132").Value = 3041.25
$ws.Range("K132").Value = 3337.7142
$ws.Range("L132").Value = 9123.75
$ws.Range("M132").Value = -807.7142000000003
$ws.Range("N132").Value = -14183.75
